$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Unnamed: 0" column (column L), shifting the columns to its
# right (M:Q) one position to the left.
$ws.Range("L1").EntireColumn.Delete()
